# Commit: "add rows for new photos"
#
# Appends 58 new filename rows (Sampo_7th_March_1.jpg .. Sampo_7th_March_58.jpg)
# to column A, right after the existing filename list, widens column A so the
# longer names fit, and leaves the selection parked at A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last populated row in column A (the filename column) and append
# the new photo rows directly below it.
$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row

for ($i = 1; $i -le 58; $i++) {
    $row = $lastRow + $i
    $ws.Cells.Item($row, 1).Value = "Sampo_7th_March_$i.jpg"
}

# Widen column A to fit the new, longer filenames.
$ws.Columns.Item(1).ColumnWidth = 21.8

# Park the selection back at A2 (top of the data, under the header row).
[void]$ws.Range("A2").Select()
